# Updates cryptos list values (price + 1h volume change) per latest scrape.
# Also fixes row order for ImmutableX / NEARProtocol (rows 37-38 swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.666.79"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "3.509.59"
$ws.Range("E3").Value = "  -2.69%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.50"
$ws.Range("E5").Value = "  -2.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.55"
$ws.Range("E6").Value = "  -3.16%  "
$ws.Range("D7").Value = "3.509.86"
$ws.Range("E7").Value = "  -2.65%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.31"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").Value = "4.110.18"
$ws.Range("E13").Value = "  -2.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.81"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("E15").Value = "  -3.88%  "
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").Value = "3.510.99"
$ws.Range("E17").Value = "  -2.56%  "
$ws.Range("D18").Value = "64.669.90"
$ws.Range("E18").Value = "  -1.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.95"
$ws.Range("E19").Value = "  -1.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.42"
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.67"
$ws.Range("E21").Value = "  -4.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.53"
$ws.Range("E22").Value = "  -1.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.578"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("D24").Value = "3.655.26"
$ws.Range("E24").Value = "  -2.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.54"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -6.93%  "
$ws.Range("E28").Value = "  -6.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.49"
$ws.Range("E29").Value = "  -8.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("E31").Value = "  -5.52%  "
$ws.Range("E32").Value = "  -5.45%  "
$ws.Range("D33").Value = "3.515.16"
$ws.Range("E33").Value = "  -2.41%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "24.08"
$ws.Range("E35").Value = "  -2.15%  "
$ws.Range("E36").Value = "  -1.37%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.60"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.26"
$ws.Range("E38").Value = "  -1.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "171.16"
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.99"
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0808"
$ws.Range("E41").Value = "  -3.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.55"
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("E43").Value = "  -3.47%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.12"
$ws.Range("E45").Value = "  -3.03%  "
$ws.Range("E46").Value = "  -3.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.39"
$ws.Range("E47").Value = "  -3.32%  "
$ws.Range("E48").Value = "  -3.01%  "
$ws.Range("D49").Value = "2.448.24"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.897"
$ws.Range("E51").Value = "  +1.42%  "
